# Update column F (dSF) values for several rows, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -4
$ws.Range("F7").Value = 4
$ws.Range("F12").Value = -10
$ws.Range("F16").Value = -4
$ws.Range("F18").Value = -5
$ws.Range("F22").Value = -2
$ws.Range("F23").Value = -6
$ws.Range("F24").Value = 1
$ws.Range("F34").Value = 7
$ws.Range("F40").Value = -4
$ws.Range("F44").Value = -3
